# Insert a new data row at row 35 (pushes existing rows 35:58 down to 36:59)
# and populate it with the new price-report record, matching the committed
# weekly update for "Hortaliza, Vega Monumental Concepción - Alcachofa".

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Rows.Item(35).Insert()

$ws.Cells.Item(35, 1).Value = 11
$ws.Cells.Item(35, 2).Value = "Vega Monumental Concepción"
$ws.Cells.Item(35, 3).Value = "Bíobío"
$ws.Cells.Item(35, 4).Value = 44777
$ws.Cells.Item(35, 5).Value = 8
$ws.Cells.Item(35, 6).Value = 100112013
$ws.Cells.Item(35, 7).Value = "Alcachofa"
$ws.Cells.Item(35, 8).Value = "Española"
$ws.Cells.Item(35, 9).Value = "Primera"
$ws.Cells.Item(35, 10).Value = 110
$ws.Cells.Item(35, 11).Value = 18000
$ws.Cells.Item(35, 12).Value = 19000
$ws.Cells.Item(35, 13).Value = 18545
$ws.Cells.Item(35, 14).Value = "$/caja 30 unidades"
$ws.Cells.Item(35, 15).Value = "Provincia de Limarí"
$ws.Cells.Item(35, 16).Value = 618
$ws.Cells.Item(35, 17).Value = 30
$ws.Cells.Item(35, 18).Value = "Hortaliza"
